$wb = $excel.ActiveWorkbook

# --- "About" sheet: add the new hydrogen-credit explanation paragraph (rows 45-48) ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A45").Value = "We also apply a very low credit for hydrogen plants because they require "
$wsAbout.Range("A46").Value = "access to a hydrogen supplier, most likely pipeline deliery, which "
$wsAbout.Range("A47").Value = "does not exist in the US today and they would only be used in certain "
$wsAbout.Range("A48").Value = "unique circumstances."

# --- "RAF-capacity" sheet: lower the regional availability factor credit for the two ---
# --- hydrogen technologies (rows 24-25, column B) from 1 to 0.5 ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")

$wsCapacity.Range("B24").Value = 0.5
$wsCapacity.Range("B25").Value = 0.5

# --- Update the on-screen selection/view state to match the saved workbook ---
$wsAbout.Activate()
$wsAbout.Range("A49").Select()

$wsCapacity.Activate()
$wsCapacity.Range("L25").Select()
